$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure values that look numeric (e.g. "1.003") are stored as text,
# matching the source data which uses inline strings for the Price/Volume columns.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "27.809.46"
$ws.Range("E2").Value = "  -0.55%  "

$ws.Range("D3").Value = "1.759.21"
$ws.Range("E3").Value = "  -1.09%  "

$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.22%  "

$ws.Range("D5").Value = "329.11"
$ws.Range("E5").Value = "  +0.68%  "

$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.14%  "

$ws.Range("D7").Value = "0.4818"
$ws.Range("E7").Value = "  +6.06%  "

$ws.Range("D8").Value = "0.3539"
$ws.Range("E8").Value = "  -1.56%  "

$ws.Range("D9").Value = "43.34"
$ws.Range("E9").Value = "  +3.19%  "

$ws.Range("D10").Value = "0.07601"
$ws.Range("E10").Value = "  +1.36%  "

$ws.Range("D11").Value = "1.077"
$ws.Range("E11").Value = "  -2.54%  "

$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  +0.11%  "

$ws.Range("D13").Value = "20.52"
$ws.Range("E13").Value = "  -1.95%  "

$ws.Range("D14").Value = "6.061"
$ws.Range("E14").Value = "  -0.01%  "

$ws.Range("D15").Value = "7.152"
$ws.Range("E15").Value = "  -0.94%  "

$ws.Range("D16").Value = "1.766.20"
$ws.Range("E16").Value = "  -0.74%  "

$ws.Range("D17").Value = "92.19"
$ws.Range("E17").Value = "  -1.79%  "

$ws.Range("D18").Value = "0.00001087"
$ws.Range("E18").Value = "  +2.04%  "

$ws.Range("D19").Value = "0.06432"
$ws.Range("E19").Value = "  -0.13%  "

$ws.Range("D20").Value = "1.002"
$ws.Range("E20").Value = "  +0.14%  "

$ws.Range("D21").Value = "16.78"
$ws.Range("E21").Value = "  -2.67%  "

$ws.Range("D22").Value = "5.778"
$ws.Range("E22").Value = "  -0.57%  "

$ws.Range("D23").Value = "27.848.26"
$ws.Range("E23").Value = "  -0.47%  "

$ws.Range("D24").Value = "11.11"
$ws.Range("E24").Value = "  -2.17%  "

$ws.Range("D25").Value = "2.153"
$ws.Range("E25").Value = "  +3.29%  "

$ws.Range("D26").Value = "163.69"
$ws.Range("E26").Value = "  -0.25%  "

$ws.Range("D27").Value = "20.03"
$ws.Range("E27").Value = "  -1.45%  "

$ws.Range("D28").Value = "1.965.07"
$ws.Range("E28").Value = "  -0.85%  "

$ws.Range("D29").Value = "2.179"
$ws.Range("E29").Value = "  -1.86%  "

$ws.Range("D30").Value = "122.90"
$ws.Range("E30").Value = "  -2.63%  "

$ws.Range("D31").Value = "1.060"
$ws.Range("E31").Value = "  -5.23%  "

$ws.Range("D32").Value = "0.09446"
$ws.Range("E32").Value = "  +2.36%  "

$ws.Range("D33").Value = "3.650"
$ws.Range("E33").Value = "  -0.84%  "

$ws.Range("D34").Value = "5.538"
$ws.Range("E34").Value = "  -0.53%  "

$ws.Range("B35").Value = "VeChain"
$ws.Range("C35").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D35").Value = "0.02257"
$ws.Range("E35").Value = "  -2.00%  "

$ws.Range("B36").Value = "Aptos"
$ws.Range("C36").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D36").Value = "11.58"
$ws.Range("E36").Value = "  -2.66%  "

$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").Value = "0.05990"
$ws.Range("E37").Value = "  -3.78%  "

$ws.Range("B38").Value = "Algorand"
$ws.Range("C38").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D38").Value = "0.2060"
$ws.Range("E38").Value = "  -1.90%  "

$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").Value = "4.892"
$ws.Range("E39").Value = "  -1.81%  "

$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Value = "0.6141"
$ws.Range("E40").Value = "  -3.42%  "

$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "1.181"
$ws.Range("E41").Value = "  -0.56%  "

$ws.Range("B42").Value = "WEMIXTOKEN"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").Value = "1.433"
$ws.Range("E42").Value = "  +3.01%  "

$ws.Range("D43").Value = "7.720"
$ws.Range("E43").Value = "  -2.66%  "

$ws.Range("D44").Value = "13.03"
$ws.Range("E44").Value = "  -2.05%  "

$ws.Range("D45").Value = "3.740"
$ws.Range("E45").Value = "  +0.04%  "

$ws.Range("D46").Value = "0.5777"
$ws.Range("E46").Value = "  -2.51%  "

$ws.Range("D47").Value = "123.54"
$ws.Range("E47").Value = "  +0.61%  "

$ws.Range("D48").Value = "1.919"
$ws.Range("E48").Value = "  -2.35%  "

$ws.Range("D49").Value = "1.153"
$ws.Range("E49").Value = "  +1.02%  "

$ws.Range("D50").Value = "0.06786"
$ws.Range("E50").Value = "  -2.11%  "

$ws.Range("D51").Value = "71.73"
$ws.Range("E51").Value = "  -1.49%  "


# Reset the explicit "text" number format back to the default style so the
# cells keep using the workbook's normal (unstyled) formatting, as before.
$ws.Range("D2:E51").Style = "Normal"
